# Update "想去人数" (want-to-go count) figures on the 展览 sheet and the
# 全部类型 rollup sheet, which mirrors the same rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5812
    $ws.Range("F3").Value = 14
    $ws.Range("F6").Value = 68
}
